$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "Бельгийский шоколад 2 плитки "
$ws.Range("B2").Value = 4250
$ws.Range("C2").Value = 45

$ws.Range("C2").Select()
